$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.898
$ws.Range("C21").Value = -12.45
$ws.Range("C23").Value = -12.68
$ws.Range("D24").Value = -7.335999999999999
$ws.Range("C25").Value = -12.063
$ws.Range("D28").Value = -8.288999999999998
$ws.Range("D36").Value = -7.417999999999999
$ws.Range("D45").Value = -7.455
$ws.Range("D48").Value = -7.540999999999999
$ws.Range("D49").Value = -8.101000000000003
$ws.Range("D52").Value = -7.733999999999999
$ws.Range("C53").Value = -11.611
$ws.Range("D53").Value = -7.567
$ws.Range("D54").Value = -7.959999999999999
$ws.Range("C57").Value = -13.742
$ws.Range("C59").Value = -12.336
$ws.Range("C69").Value = -10.921
$ws.Range("D70").Value = -7.360000000000001
$ws.Range("C79").Value = -12.338
$ws.Range("C83").Value = -13.351
$ws.Range("D86").Value = -8.294
$ws.Range("D87").Value = -8.053000000000001
$ws.Range("C93").Value = -10.836
$ws.Range("D101").Value = -7.233
